$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 18 de Septiembre de 2020 a las 02:31"

# Row 4
$ws.Range("B4").Value = 6874023
$ws.Range("C4").Value = 45722
$ws.Range("D4").Value = 4150492
$ws.Range("E4").Value = 2521341
$ws.Range("G4").Value = 856
$ws.Range("H4").Value = 202190

# Row 25
$ws.Range("B25").Value = 269042
$ws.Range("C25").Value = 2177
$ws.Range("D25").Value = 241300
$ws.Range("E25").Value = 18285
$ws.Range("G25").Value = 8
$ws.Range("H25").Value = 9457

# Row 66
$ws.Range("B66").Value = 44155
$ws.Range("C66").Value = 3123
$ws.Range("D66").Value = 23377
$ws.Range("E66").Value = 20289
$ws.Range("G66").Value = 7
$ws.Range("H66").Value = 489

# Row 107
$ws.Range("B107").Value = 7541
$ws.Range("C107").Value = 147
$ws.Range("D107").Value = 6641
$ws.Range("E107").Value = 776

# Row 108
$ws.Range("A108").Value = "Montenegro"
$ws.Range("B108").Value = 7503
$ws.Range("C108").Value = 212
$ws.Range("D108").Value = 4892
$ws.Range("E108").Value = 2482
$ws.Range("G108").Value = 3
$ws.Range("H108").Value = 129

# Row 109
$ws.Range("A109").Value = "Mauritania"
$ws.Range("B109").Value = 7346
$ws.Range("C109").Value = 14
$ws.Range("D109").Value = 6865
$ws.Range("E109").Value = 320
$ws.Range("H109").Value = 161

# Row 115
$ws.Range("B115").Value = 5191
$ws.Range("C115").Value = 36
$ws.Range("D115").Value = 4439
$ws.Range("E115").Value = 649
$ws.Range("G115").Value = 2
$ws.Range("H115").Value = 103

# Row 123
$ws.Range("A123").Value = "Surinam"
$ws.Range("B123").Value = 4671
$ws.Range("C123").Value = 26
$ws.Range("D123").Value = 4160
$ws.Range("E123").Value = 416
$ws.Range("G123").Value = 0
$ws.Range("H123").Value = 95

# Row 124
$ws.Range("A124").Value = "Ruanda"
$ws.Range("B124").Value = 4653
$ws.Range("C124").Value = 19
$ws.Range("D124").Value = 2817
$ws.Range("E124").Value = 1813
$ws.Range("G124").Value = 1
$ws.Range("H124").Value = 23

# Row 131
$ws.Range("A131").Value = "Mayotte"
$ws.Range("B131").Value = 3541
$ws.Range("C131").Value = 0
$ws.Range("D131").Value = 2964
$ws.Range("E131").Value = 537
$ws.Range("H131").Value = 40

# Row 132
$ws.Range("A132").Value = "Lituania"
$ws.Range("B132").Value = 3504
$ws.Range("C132").Value = 62
$ws.Range("D132").Value = 2149
$ws.Range("E132").Value = 1268
$ws.Range("H132").Value = 87

# Row 133
$ws.Range("A133").Value = "Tailandia"
$ws.Range("B133").Value = 3490
$ws.Range("D133").Value = 3325
$ws.Range("E133").Value = 107
$ws.Range("H133").Value = 58

# Row 134
$ws.Range("A134").Value = "Gambia"
$ws.Range("B134").Value = 3473
$ws.Range("C134").Value = 33
$ws.Range("D134").Value = 1951
$ws.Range("E134").Value = 1415
$ws.Range("G134").Value = 0
$ws.Range("H134").Value = 107

# Row 135
$ws.Range("A135").Value = "Trinidad yTobago"
$ws.Range("B135").Value = 3434
$ws.Range("C135").Value = 107
$ws.Range("D135").Value = 1469
$ws.Range("E135").Value = 1905
$ws.Range("G135").Value = 2
$ws.Range("H135").Value = 60

# Row 136
$ws.Range("A136").Value = "Guadalupe"
$ws.Range("B136").Value = 3426
$ws.Range("D136").Value = 837
$ws.Range("E136").Value = 2563
$ws.Range("H136").Value = 26

# Row 137
$ws.Range("A137").Value = "Somalia"
$ws.Range("B137").Value = 3390
$ws.Range("D137").Value = 2812
$ws.Range("E137").Value = 480
$ws.Range("H137").Value = 98

# Row 138
$ws.Range("B138").Value = 3382
$ws.Range("C138").Value = 54
$ws.Range("D138").Value = 1742
$ws.Range("E138").Value = 1617
$ws.Range("G138").Value = 1
$ws.Range("H138").Value = 23

# Row 139
$ws.Range("B139").Value = 3276
$ws.Range("C139").Value = 5
$ws.Range("E139").Value = 220

# Row 140
$ws.Range("A140").Value = "Bahamas"
$ws.Range("B140").Value = 3177
$ws.Range("C140").Value = 90
$ws.Range("D140").Value = 1626
$ws.Range("E140").Value = 1482
$ws.Range("H140").Value = 69

# Row 141
$ws.Range("A141").Value = "Reunion"
$ws.Range("B141").Value = 3099
$ws.Range("C141").Value = 97
$ws.Range("D141").Value = 1794
$ws.Range("E141").Value = 1290
$ws.Range("H141").Value = 15

# Row 145
$ws.Range("A145").Value = "Sudan del Sur"
$ws.Range("B145").Value = 2599
$ws.Range("C145").Value = 5
$ws.Range("D145").Value = 1290
$ws.Range("E145").Value = 1260
$ws.Range("H145").Value = 49

# Row 146
$ws.Range("A146").Value = "Malta"
$ws.Range("B146").Value = 2595
$ws.Range("C146").Value = 35
$ws.Range("D146").Value = 1978
$ws.Range("E146").Value = 601
$ws.Range("H146").Value = 16

# Row 154
$ws.Range("B154").Value = 1876
$ws.Range("C154").Value = 20
$ws.Range("D154").Value = 1582
$ws.Range("E154").Value = 248
$ws.Range("G154").Value = 1
$ws.Range("H154").Value = 46

# Row 169
$ws.Range("D169").Value = 873
$ws.Range("E169").Value = 19

# Row 203
$ws.Range("D203").Value = 26
$ws.Range("E203").Value = 4

# Row 214
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0

# Row 215
$ws.Range("A215").Value = "Montserrat"
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1
